$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Differences" columns (N/O) added alongside rows 11-13 ---
$ws.Range("N11").Formula = "=D4 - E4"
$ws.Range("O11").Formula = "= N11 * O4"

$ws.Range("N12").Formula = "=D5 - E5"
$ws.Range("O12").Formula = "= N12 * O5"

$ws.Range("N13").Formula = "=D6 - E6"
$ws.Range("O13").Formula = "= N13 * O6"

# --- Remove the stray "f" label that used to sit in O21 ---
$ws.Range("O21").ClearContents()

# --- New "Averages" block starting at row 33 ---
$ws.Range("A33").Value = "Averages"

$ws.Range("A34").Value = "5.46m"
$ws.Range("B34").Value = "ws"
$ws.Range("D34").Value = "u"
$ws.Range("E34").Value = "v"

$ws.Range("B35").Formula = "=AVERAGE(D2:D4) + 3"
$ws.Range("D35").Value = 7.702
$ws.Range("E35").Value = -9.3960000000000008
$ws.Range("F35").Value = 12.1496116
$ws.Range("G35").Formula = "=SQRT((POWER(D35,2)+POWER(E35,2)))"
$ws.Range("H35").Formula = "=AVERAGE(G18:G20)"
$ws.Range("I35").Formula = "=AVERAGE(H18:H20)"
$ws.Range("L35").Value = 289.07
$ws.Range("M35").Formula = "=AVERAGE(L35:L37)"

$ws.Range("B36").Formula = "=D5 + 3"
$ws.Range("D36").Value = 8.1199999999999992
$ws.Range("E36").Value = -10.917
$ws.Range("F36").Value = 13.6060832
$ws.Range("G36").Formula = "=SQRT((POWER(D36,2)+POWER(E36,2)))"
$ws.Range("H36").Value = 53.356993750000001
$ws.Range("H36").Font.Color = 0
$ws.Range("I36").Value = 36.643006200000002
$ws.Range("I36").Font.Color = 0
$ws.Range("L36").Value = 288.23

$ws.Range("B37").Formula = "=D6 + O13"
$ws.Range("D37").Value = 9.9700000000000006
$ws.Range("E37").Value = -13.78
$ws.Range("G37").Formula = "=SQRT((POWER(D37,2)+POWER(E37,2)))"
$ws.Range("H37").Value = 54.121180070000001
$ws.Range("H37").Font.Color = 0
$ws.Range("I37").Value = 35.878819900000003
$ws.Range("I37").Font.Color = 0
$ws.Range("L37").Value = 287.05

# --- Selection / view bookkeeping to mirror the saved workbook state ---
$ws.Range("F36").Select()
